$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the old row 7 (shifts the two header rows
# down by one: old row 8 -> 9, old row 9 -> 10) and leaves room for a new
# data row at 11.
$ws.Rows(7).Insert()

# Pre-format the date cells with the built-in short-date number format
# BEFORE writing values into them, so Excel reuses a single style (numFmt
# 14) instead of inventing a fresh one per cell.
$ws.Range("B11,E11,T11,U11").NumberFormat = "mm-dd-yy"

# New data row (row 11)
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "2018-05-01"
$ws.Range("C11").Value = "ext. ID 123"
$ws.Range("D11").Value = "LAB123"
$ws.Range("E11").Value = "2018-01-01"
$ws.Range("F11").Value = 10
$ws.Range("G11").Value = "ko/ko"
$ws.Range("P11").Value = "m"
$ws.Range("Q11").Value = "Outside"
$ws.Range("R11").Value = "#007"
$ws.Range("S11").Value = "Roger Roger"
$ws.Range("T11").Value = "2018-06-01"
$ws.Range("U11").Value = "2018-07-01"

# Update the view so the active selection matches the saved file
# (scrolled so column C is left-most, cell F12 selected).
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("F12").Select()
